$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Edit 1: Rework the "Robot Framework" list-paragraph block.
#   - moves "Test Setup & Tear Down..." text to share the lastRenderedPageBreak
#     run and relocates the remainder of the block, adding the new
#     "Test Template" / Data-Driven-via-Excel sub bullets and example block.
# ---------------------------------------------------------------------------
$anchor = $d.Content
$found = $anchor.Find.Execute("Test Setup & Tear Down")
if (-not $found) {
    throw "Could not find 'Test Setup & Tear Down' anchor text"
}
$anchorPara = $anchor.Paragraphs(1)
$startPara = $anchorPara.Previous(2)
$endPara = $anchorPara.Next()
$targetRange = $d.Range($startPara.Range.Start, $endPara.Range.End)

$newRegionXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="4"/></w:numPr><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:lastRenderedPageBreak/><w:t xml:space="preserve"> Test Setup &amp; Tear Down – runs before and after each test case</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="4"/></w:numPr><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t xml:space="preserve">Robot Framework </w:t></w:r><w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:tab/></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="4"/></w:numPr><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t xml:space="preserve">Data-Driven Framework </w:t></w:r><w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t xml:space="preserve">– Test Template </w:t></w:r><w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>in robot framework</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="4"/></w:numPr><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t xml:space="preserve">Create a keyword with test data as an arguments </w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="4"/></w:numPr><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t xml:space="preserve">Declare the template in setting sections </w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="4"/></w:numPr><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t xml:space="preserve">Create Test case and pass the test data(arguments) </w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="4"/></w:numPr><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t xml:space="preserve">Data-Driven using excel – Test Template &amp; Excel </w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="4"/></w:numPr><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t xml:space="preserve">Install </w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:ind w:left="2070"/><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas"/><w:color w:val="0E1116"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas"/><w:color w:val="0E1116"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>pip install --upgrade robotframework-datadriver</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:ind w:left="2070"/><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas"/><w:color w:val="0E1116"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas"/><w:color w:val="0E1116"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>pip install --upgrade robotframework-datadriver[XLS]</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="4"/></w:numPr><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas"/><w:color w:val="0E1116"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>Make sure the arguments in test template in present as a header in the excel sheet</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="4"/></w:numPr><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t xml:space="preserve">Declare settings section </w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="HTMLPreformatted"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/><w:ind w:left="720"/><w:rPr><w:color w:val="080808"/><w:sz w:val="30"/><w:szCs w:val="30"/></w:rPr></w:pPr><w:r><w:rPr><w:color w:val="9E880D"/><w:sz w:val="30"/><w:szCs w:val="30"/></w:rPr><w:tab/></w:r><w:r><w:rPr><w:color w:val="9E880D"/><w:sz w:val="30"/><w:szCs w:val="30"/></w:rPr><w:tab/></w:r><w:bookmarkStart w:id="3" w:name="_GoBack"/><w:bookmarkEnd w:id="3"/><w:r><w:rPr><w:color w:val="9E880D"/><w:sz w:val="30"/><w:szCs w:val="30"/></w:rPr><w:t xml:space="preserve">Library     </w:t></w:r><w:r><w:rPr><w:color w:val="067D17"/><w:sz w:val="30"/><w:szCs w:val="30"/></w:rPr><w:t xml:space="preserve">DataDriver      </w:t></w:r><w:r><w:rPr><w:color w:val="1750EB"/><w:sz w:val="30"/><w:szCs w:val="30"/></w:rPr><w:t>file</w:t></w:r><w:r><w:rPr><w:color w:val="080808"/><w:sz w:val="30"/><w:szCs w:val="30"/></w:rPr><w:t xml:space="preserve">=../test_data/openemr_data.xlsx     </w:t></w:r><w:r><w:rPr><w:color w:val="1750EB"/><w:sz w:val="30"/><w:szCs w:val="30"/></w:rPr><w:t>sheet_name</w:t></w:r><w:r><w:rPr><w:color w:val="080808"/><w:sz w:val="30"/><w:szCs w:val="30"/></w:rPr><w:t>=InvalidLoginTest</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:ind w:left="2070"/><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:ind w:left="2070"/><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:ind w:left="1350"/><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr></w:pPr></w:p>
'@
$targetRange.InsertXML($newRegionXml)

# ---------------------------------------------------------------------------
# Edit 2: add back a <w:lastRenderedPageBreak/> in the "partial link" table
# cell run (highlighted green) which now receives the page break that moved.
# ---------------------------------------------------------------------------
$linkRange = $d.Content
$found2 = $linkRange.Find.Execute("partial link")
if (-not $found2) {
    throw "Could not find 'partial link' anchor text"
}
$linkPara = $linkRange.Paragraphs(1)
$linkRunRange = $linkPara.Range
$linkXml = @'
<w:r xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/><w:color w:val="000000"/><w:highlight w:val="green"/></w:rPr><w:lastRenderedPageBreak/><w:t>partial link</w:t></w:r>
'@
$linkRunRange.InsertXML($linkXml)

# ---------------------------------------------------------------------------
# Edit 3: remove the <w:lastRenderedPageBreak/> that used to precede
# "Select Job title as ..." (the break shifted earlier in the document).
# ---------------------------------------------------------------------------
$jobRange = $d.Content
$found3 = $jobRange.Find.Execute("Select Job title as")
if (-not $found3) {
    throw "Could not find 'Select Job title as' anchor text"
}
$jobPara = $jobRange.Paragraphs(1)
$jobXml = @'
<w:r xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:rPr><w:color w:val="000000"/></w:rPr><w:t>Select Job title as “IT Manager”</w:t></w:r>
'@
$jobPara.Range.InsertXML($jobXml)

Write-Host "Done"
